$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.33"
$ws.Range("E2").Value = "'-0.77%"
$ws.Range("D3").Value = "'31.07"
$ws.Range("E3").Value = "'0.86%"
$ws.Range("D4").Value = "'4.927"
$ws.Range("E4").Value = "'-0.64%"
$ws.Range("D5").Value = "'0.07352"
$ws.Range("E5").Value = "'1.81%"
$ws.Range("D6").Value = "'2.254"
$ws.Range("E6").Value = "'25.76%"
$ws.Range("D7").Value = "'7.724"
$ws.Range("E7").Value = "'0.57%"
$ws.Range("D8").Value = "'3.725"
$ws.Range("E8").Value = "'-0.97%"
$ws.Range("D9").Value = "'0.9049"
$ws.Range("E9").Value = "'0.96%"
$ws.Range("D10").Value = "'0.08935"
$ws.Range("E10").Value = "'16.04%"
$ws.Range("D11").Value = "'0.1689"
$ws.Range("E11").Value = "'2.38%"
$ws.Range("D12").Value = "'0.08178"
$ws.Range("E12").Value = "'2.50%"
$ws.Range("D13").Value = "'0.03114"
$ws.Range("E13").Value = "'2.40%"
$ws.Range("D14").Value = "'0.09947"
$ws.Range("E14").Value = "'-0.83%"
$ws.Range("E15").Value = "'-0.39%"
$ws.Range("D16").Value = "'0.005780"
$ws.Range("E16").Value = "'0.41%"
$ws.Range("D17").Value = "'3.485"
$ws.Range("E17").Value = "'0.50%"
$ws.Range("D18").Value = "'2.071"
$ws.Range("E18").Value = "'-0.54%"
$ws.Range("E19").Value = "'0.53%"
$ws.Range("D20").Value = "'0.1288"
$ws.Range("E20").Value = "'-1.04%"
$ws.Range("D21").Value = "'4.159"
$ws.Range("E21").Value = "'3.10%"
$ws.Range("E22").Value = "'-9.55%"
$ws.Range("D23").Value = "'0.04550"
$ws.Range("E23").Value = "'0.90%"
$ws.Range("D24").Value = "'0.001208"
$ws.Range("E24").Value = "'-0.47%"
$ws.Range("D25").Value = "'0.004157"
$ws.Range("E25").Value = "'3.56%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'3.93%"
$ws.Range("D39").Value = "'0.01582"
$ws.Range("E39").Value = "'-1.33%"
$ws.Range("D40").Value = "'0.04452"
$ws.Range("E40").Value = "'1.30%"
$ws.Range("D41").Value = "'0.007361"
$ws.Range("E41").Value = "'0.98%"
$ws.Range("D42").Value = "'0.009655"
$ws.Range("E42").Value = "'25.75%"
$ws.Range("E43").Value = "'1.30%"
$ws.Range("E44").Value = "'17.29%"
$ws.Range("D45").Value = "'0.008514"
$ws.Range("E45").Value = "'-7.52%"
$ws.Range("D46").Value = "'0.00006113"
$ws.Range("E46").Value = "'2.91%"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("D48").Value = "'2.245"
$ws.Range("E48").Value = "'-0.04%"
$ws.Range("E49").Value = "'-33.38%"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("E51").Value = "'-0.05%"
